# Rename "scene" to "content" in the Binary glTF extension figure slide.
#
# Three labelled boxes on slide 1 reference the old "scene*" field names
# from the glTF container header diagram; rename them to the new
# "content*" names while preserving all existing run-level formatting
# (font, color, size, etc.).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "CustomShape 3": sceneFormat(uint32)  ->  contentFormat(uint32)
$shapeFormat = $s.Shapes.Item(3)
$shapeFormat.TextFrame.TextRange.Replace("sceneFormat", "contentFormat")

# "CustomShape 4": sceneLength(uint32)  ->  contentLength(uint32)
# Only the "scene" prefix is renamed to "content"; the "Length" suffix is
# left as-is but ends up as its own run because the replaced prefix has a
# different length, matching how PowerPoint splits runs on a partial edit.
$shapeLength = $s.Shapes.Item(4)
$sceneChars = $shapeLength.TextFrame.TextRange.Characters(1, 5)
$sceneChars.Text = "content"

# "CustomShape 14": scene  ->  content
$shapeScene = $s.Shapes.Item(14)
$shapeScene.TextFrame.TextRange.Replace("scene", "content")
